$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("H98").Value = 1625.875
$ws.Range("J98").Value = 2950.2
$ws.Range("L98").Value = 2950.2
$ws.Range("N98").Value = -5946.2
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H108").Value = 99999
$ws.Range("J108").Value = 99999
$ws.Range("L108").Value = 99999
$ws.Range("N108").Value = -107679
$ws.Range("H109").Value = 99999
$ws.Range("J109").Value = 99999
$ws.Range("L109").Value = 99999
$ws.Range("N109").Value = -102773
$ws.Range("H112").Value = 1289.5
$ws.Range("J112").Value = 1639
$ws.Range("L112").Value = 4917
$ws.Range("N112").Value = -7133
$ws.Range("H117").Value = 94840.75
$ws.Range("J117").Value = 94840.75
$ws.Range("L117").Value = 94840.75
$ws.Range("N117").Value = -104018.75
$ws.Range("H122").Value = 1625.875
$ws.Range("J122").Value = 2950.2
$ws.Range("L122").Value = 8850.599999999999
$ws.Range("N122").Value = -13750.6
$ws.Range("H123").Value = 70737.5
$ws.Range("J123").Value = 70737.5
$ws.Range("L123").Value = 70737.5
$ws.Range("N123").Value = -80537.5
$ws.Range("H133").Value = 69517.39999999999
$ws.Range("J133").Value = 69517.39999999999
$ws.Range("L133").Value = 69517.39999999999
$ws.Range("N133").Value = -79637.39999999999
$ws.Range("H136").Value = 85324.336
$ws.Range("J136").Value = 85324.336
$ws.Range("L136").Value = 85324.336
$ws.Range("N136").Value = -95524.336
$ws.Range("H139").Value = 98402
$ws.Range("J139").Value = 98402
$ws.Range("L139").Value = 98402
$ws.Range("N139").Value = -108682
$ws.Range("H140").Value = 80776
$ws.Range("J140").Value = 80776
$ws.Range("L140").Value = 80776
$ws.Range("N140").Value = -91136
$ws.Range("M8").ClearContents()
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 52620.89
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 52620.89
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 52620.89
$ws.Range("N52").Value = -53256.89
$ws.Range("H110").Value = 1107.7059
$ws.Range("I110").Value = 1054.7693
$ws.Range("J110").Value = 1279.75
$ws.Range("K110").Value = 1054.7693
$ws.Range("L110").Value = 1279.75
$ws.Range("M110").Value = 990.2307000000001
$ws.Range("N110").Value = -5369.75
$ws.Range("H117").Value = 42730
$ws.Range("J117").Value = 42730
$ws.Range("L117").Value = 42730
$ws.Range("N117").Value = -51908
$ws.Range("H121").Value = 48960
$ws.Range("J121").Value = 48960
$ws.Range("L121").Value = 48960
$ws.Range("N121").Value = -52454
$ws.Range("M52").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 34665.332
$ws.Range("J55").Value = 34665.332
$ws.Range("L55").Value = 34665.332
$ws.Range("N55").Value = -35211.332
$ws.Range("H108").Value = 94621.125
$ws.Range("J108").Value = 94621.125
$ws.Range("L108").Value = 94621.125
$ws.Range("N108").Value = -102301.125
$ws.Range("H110").Value = 83354.28999999999
$ws.Range("J110").Value = 83354.28999999999
$ws.Range("L110").Value = 83354.28999999999
$ws.Range("N110").Value = -91534.28999999999
$ws.Range("H132").Value = 28535.5
$ws.Range("J132").Value = 28535.5
$ws.Range("L132").Value = 28535.5
$ws.Range("N132").Value = -38655.5
$ws.Range("H135").Value = 118696
$ws.Range("J135").Value = 118696
$ws.Range("L135").Value = 118696
$ws.Range("N135").Value = -128836
$ws.Range("H138").Value = 99758
$ws.Range("J138").Value = 99758
$ws.Range("L138").Value = 99758
$ws.Range("N138").Value = -110038

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 32298.4
$ws.Range("J9").Value = 32298.4
$ws.Range("L9").Value = 32298.4
$ws.Range("N9").Value = -32634.4
$ws.Range("H18").Value = 24226
$ws.Range("J18").Value = 24226
$ws.Range("L18").Value = 24226
$ws.Range("N18").Value = -24686
$ws.Range("H31").Value = 4040.6667
$ws.Range("I31").Value = 2992.8572
$ws.Range("J31").Value = 5507.6
$ws.Range("K31").Value = 2992.8572
$ws.Range("L31").Value = 5507.6
$ws.Range("M31").Value = -2697.8572
$ws.Range("N31").Value = -6097.6
$ws.Range("H34").Value = 4040.6667
$ws.Range("I34").Value = 2992.8572
$ws.Range("J34").Value = 5507.6
$ws.Range("K34").Value = 2992.8572
$ws.Range("L34").Value = 5507.6
$ws.Range("M34").Value = -2790.8572
$ws.Range("N34").Value = -5911.6
$ws.Range("H86").Value = 8765
$ws.Range("I86").Value = 2995
$ws.Range("J86").Value = 11650
$ws.Range("K86").Value = 2995
$ws.Range("L86").Value = 11650
$ws.Range("M86").Value = -1872
$ws.Range("N86").Value = -13896
$ws.Range("H89").Value = 8765
$ws.Range("I89").Value = 2995
$ws.Range("J89").Value = 11650
$ws.Range("K89").Value = 14975
$ws.Range("L89").Value = 58250
$ws.Range("M89").Value = -9359
$ws.Range("N89").Value = -69482
$ws.Range("H108").Value = 57395.363
$ws.Range("J108").Value = 57395.363
$ws.Range("L108").Value = 57395.363
$ws.Range("N108").Value = -65075.363
$ws.Range("H117").Value = 37641.145
$ws.Range("J117").Value = 37641.145
$ws.Range("L117").Value = 37641.145
$ws.Range("N117").Value = -46819.145
$ws.Range("H122").Value = 2249.7856
$ws.Range("I122").Value = 1949.95
$ws.Range("J122").Value = 2999.375
$ws.Range("K122").Value = 5849.85
$ws.Range("L122").Value = 8998.125
$ws.Range("M122").Value = -3399.85
$ws.Range("N122").Value = -13898.125
$ws.Range("H132").Value = 2338.238
$ws.Range("I132").Value = 2036.9333
$ws.Range("K132").Value = 6110.7999
$ws.Range("M132").Value = -3580.7999
$ws.Range("H138").Value = 94492
$ws.Range("J138").Value = 94492
$ws.Range("L138").Value = 94492
$ws.Range("N138").Value = -104772

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1758.2
$ws.Range("I36").Value = 2317.3333
$ws.Range("J36").Value = 919.5
$ws.Range("K36").Value = 6951.999899999999
$ws.Range("L36").Value = 2758.5
$ws.Range("M36").Value = -6782.999899999999
$ws.Range("N36").Value = -3096.5
$ws.Range("H60").Value = 1373.0555
$ws.Range("I60").Value = 190
$ws.Range("J60").Value = 1609.6666
$ws.Range("K60").Value = 570
$ws.Range("L60").Value = 4828.9998
$ws.Range("M60").Value = -319
$ws.Range("N60").Value = -5330.9998
$ws.Range("H87").Value = 11793.5
$ws.Range("I87").Value = 1087.5
$ws.Range("K87").Value = 3262.5
$ws.Range("M87").Value = -2014.5
$ws.Range("H90").Value = 11793.5
$ws.Range("I90").Value = 1087.5
$ws.Range("K90").Value = 9787.5
$ws.Range("M90").Value = -3547.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20311.7
$ws.Range("J93").Value = 20311.7
$ws.Range("L93").Value = 20311.7
$ws.Range("N93").Value = -24055.7
$ws.Range("H108").Value = 51801.168
$ws.Range("J108").Value = 51801.168
$ws.Range("L108").Value = 51801.168
$ws.Range("N108").Value = -59481.168
$ws.Range("H114").Value = 63463.363
$ws.Range("J114").Value = 63463.363
$ws.Range("L114").Value = 63463.363
$ws.Range("N114").Value = -72141.363
$ws.Range("H132").Value = 3752.5557
$ws.Range("I132").Value = 2945.5557
$ws.Range("K132").Value = 8836.667099999999
$ws.Range("M132").Value = -6306.667099999999
$ws.Range("H140").Value = 90430.664
$ws.Range("J140").Value = 90396
$ws.Range("L140").Value = 90396
$ws.Range("N140").Value = -100756

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 73278
$ws.Range("J129").Value = 72500
$ws.Range("L129").Value = 72500
$ws.Range("N129").Value = -82500

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 115384.54
$ws.Range("I24").Value = 1000000
$ws.Range("J24").Value = 41666.582
$ws.Range("K24").Value = 1000000
$ws.Range("L24").Value = 41666.582
$ws.Range("M24").Value = -999770
$ws.Range("N24").Value = -42126.582
$ws.Range("H121").Value = 36997.332
$ws.Range("J121").Value = 36997.332
$ws.Range("L121").Value = 36997.332
$ws.Range("N121").Value = -40491.332
$ws.Range("H127").Value = 87177.14
$ws.Range("J127").Value = 91641.664
$ws.Range("L127").Value = 91641.664
$ws.Range("N127").Value = -101561.664
